# Update "想去人数" (F column) values on the "展览" and "全部类型" sheets
# F3: 183 -> 184
# F4: 135 -> 137

$wb = $excel.ActiveWorkbook

$sheetNames = @("展览", "全部类型")

foreach ($name in $sheetNames) {
    $ws = $wb.Worksheets.Item($name)
    $ws.Range("F3").Value = 184
    $ws.Range("F4").Value = 137
}
